$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently in row 535 (columns that remain unchanged after the
# new row is inserted: Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria,
# Variedad, Calidad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion)
$valA = $ws.Cells.Item(535, 1).Value2
$valB = $ws.Cells.Item(535, 2).Value2
$valC = $ws.Cells.Item(535, 3).Value2
$valE = $ws.Cells.Item(535, 5).Value2
$valF = $ws.Cells.Item(535, 6).Value2
$valG = $ws.Cells.Item(535, 7).Value2
$valH = $ws.Cells.Item(535, 8).Value2
$valI = $ws.Cells.Item(535, 9).Value2
$valN = $ws.Cells.Item(535, 14).Value2
$valO = $ws.Cells.Item(535, 15).Value2
$valQ = $ws.Cells.Item(535, 17).Value2
$valR = $ws.Cells.Item(535, 18).Value2

# Insert a new row at 535, pushing the existing rows 535:624 down to 536:625
$ws.Rows.Item(535).Insert()

# Rebuild row 535 with the carried-over values plus the new weekly record values
$ws.Cells.Item(535, 1).Value = $valA
$ws.Cells.Item(535, 2).Value = $valB
$ws.Cells.Item(535, 3).Value = $valC
$ws.Cells.Item(535, 4).Value = 45180
$ws.Cells.Item(535, 5).Value = $valE
$ws.Cells.Item(535, 6).Value = $valF
$ws.Cells.Item(535, 7).Value = $valG
$ws.Cells.Item(535, 8).Value = $valH
$ws.Cells.Item(535, 9).Value = $valI
$ws.Cells.Item(535, 10).Value = 420
$ws.Cells.Item(535, 11).Value = 15000
$ws.Cells.Item(535, 12).Value = 16000
$ws.Cells.Item(535, 13).Value = 15429
$ws.Cells.Item(535, 14).Value = $valN
$ws.Cells.Item(535, 15).Value = $valO
$ws.Cells.Item(535, 16).Value = 309
$ws.Cells.Item(535, 17).Value = $valQ
$ws.Cells.Item(535, 18).Value = $valR
